# mayers_mike.xlsx regen: replace column G (header "K") values that were
# previously populated from "Strike#" with the recomputed strikeout counts.
# Only column G (rows 2-83, skipping rows whose value is unchanged) differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K value, taken from the regenerated save_data.
$kValues = [ordered]@{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 3
    9  = 2
    10 = 3
    11 = 2
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 4
    20 = 0
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 2
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    36 = 1
    37 = 3
    38 = 2
    39 = 3
    40 = 2
    41 = 3
    42 = 0
    43 = 3
    44 = 1
    45 = 0
    46 = 2
    48 = 1
    49 = 1
    50 = 1
    51 = 1
    52 = 2
    53 = 1
    54 = 3
    55 = 3
    56 = 1
    57 = 3
    58 = 1
    59 = 0
    60 = 2
    61 = 1
    62 = 1
    63 = 2
    64 = 1
    65 = 0
    66 = 4
    67 = 0
    68 = 4
    69 = 1
    70 = 1
    71 = 0
    72 = 2
    73 = 1
    74 = 2
    75 = 3
    76 = 1
    77 = 1
    78 = 1
    79 = 1
    80 = 1
    81 = 2
    83 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
